# Scheduled runner update: refresh computed market-board profit figures
# (currentAveragePrice / Price / Profit columns H:N) across several
# crafting-leve sheets, matching the latest price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 15629084
$ws.Range("I132").Value = 17548070
$ws.Range("J132").Value = 3063.2856
$ws.Range("K132").Value = 52644210
$ws.Range("L132").Value = 9189.856800000001
$ws.Range("M132").Value = -52641680
$ws.Range("N132").Value = -14249.8568
$ws.Range("H134").Value = 103690.266
$ws.Range("J134").Value = 103690.266
$ws.Range("L134").Value = 103690.266
$ws.Range("N134").Value = -113830.266
$ws.Range("H135").Value = 564.03705
$ws.Range("I135").Value = 419.54166
$ws.Range("J135").Value = 1720
$ws.Range("K135").Value = 3775.87494
$ws.Range("L135").Value = 15480
$ws.Range("M135").Value = -1240.87494
$ws.Range("N135").Value = -20550
$ws.Range("H137").Value = 50104.676
$ws.Range("I137").Value = 61274.066
$ws.Range("J137").Value = 2235.8572
$ws.Range("K137").Value = 183822.198
$ws.Range("L137").Value = 6707.571599999999
$ws.Range("M137").Value = -181272.198
$ws.Range("N137").Value = -11807.5716
$ws.Range("H138").Value = 2894.215
$ws.Range("I138").Value = 1019.5172
$ws.Range("J138").Value = 3981.54
$ws.Range("K138").Value = 3058.5516
$ws.Range("L138").Value = 11944.62
$ws.Range("M138").Value = 2081.4484
$ws.Range("N138").Value = -22224.62

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3708270
$ws.Range("I94").Value = 4167969.5
$ws.Range("J94").Value = 30673.334
$ws.Range("K94").Value = 4167969.5
$ws.Range("L94").Value = 30673.334
$ws.Range("M94").Value = -4167518.5
$ws.Range("N94").Value = -31575.334
$ws.Range("H105").Value = 1786920.9
$ws.Range("I105").Value = 2017062.9
$ws.Range("J105").Value = 3321
$ws.Range("K105").Value = 2017062.9
$ws.Range("L105").Value = 3321
$ws.Range("M105").Value = -2015315.9
$ws.Range("N105").Value = -6815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19452.176
$ws.Range("I31").Value = 3105.3438
$ws.Range("J31").Value = 40376.12
$ws.Range("K31").Value = 3105.3438
$ws.Range("L31").Value = 40376.12
$ws.Range("M31").Value = -2810.3438
$ws.Range("N31").Value = -40966.12
$ws.Range("H34").Value = 19452.176
$ws.Range("I34").Value = 3105.3438
$ws.Range("J34").Value = 40376.12
$ws.Range("K34").Value = 3105.3438
$ws.Range("L34").Value = 40376.12
$ws.Range("M34").Value = -2903.3438
$ws.Range("N34").Value = -40780.12
$ws.Range("H132").Value = 56441.133
$ws.Range("I132").Value = 32370.25
$ws.Range("J132").Value = 184819.17
$ws.Range("K132").Value = 97110.75
$ws.Range("L132").Value = 554457.51
$ws.Range("M132").Value = -94580.75
$ws.Range("N132").Value = -559517.51
$ws.Range("H134").Value = 17761.053
$ws.Range("I134").Value = 23095.365
$ws.Range("J134").Value = 3180.6
$ws.Range("K134").Value = 69286.095
$ws.Range("L134").Value = 9541.799999999999
$ws.Range("M134").Value = -66751.095
$ws.Range("N134").Value = -14611.8
$ws.Range("H141").Value = 43293.3
$ws.Range("J141").Value = 43293.3
$ws.Range("L141").Value = 43293.3
$ws.Range("N141").Value = -53653.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1985079.4
$ws.Range("I97").Value = 5953358
$ws.Range("J97").Value = 940
$ws.Range("K97").Value = 5953358
$ws.Range("L97").Value = 940
$ws.Range("M97").Value = -5952862
$ws.Range("N97").Value = -1932

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2028.3572
$ws.Range("I16").Value = 850
$ws.Range("J16").Value = 2499.7
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 2499.7
$ws.Range("M16").Value = -680
$ws.Range("N16").Value = -2839.7
$ws.Range("H22").Value = 37803.457
$ws.Range("I22").Value = 111958.25
$ws.Range("J22").Value = 726.0625
$ws.Range("K22").Value = 111958.25
$ws.Range("L22").Value = 726.0625
$ws.Range("M22").Value = -111663.25
$ws.Range("N22").Value = -1316.0625
$ws.Range("H27").Value = 37803.457
$ws.Range("I27").Value = 111958.25
$ws.Range("J27").Value = 726.0625
$ws.Range("K27").Value = 111958.25
$ws.Range("L27").Value = 726.0625
$ws.Range("M27").Value = -111851.25
$ws.Range("N27").Value = -940.0625
$ws.Range("H40").Value = 6816.6665
$ws.Range("I40").Value = 3982
$ws.Range("J40").Value = 9651.333000000001
$ws.Range("K40").Value = 3982
$ws.Range("L40").Value = 9651.333000000001
$ws.Range("M40").Value = -3846
$ws.Range("N40").Value = -9923.333000000001
$ws.Range("H46").Value = 6785.727
$ws.Range("I46").Value = 4940.2
$ws.Range("J46").Value = 8323.666999999999
$ws.Range("K46").Value = 4940.2
$ws.Range("L46").Value = 8323.666999999999
$ws.Range("M46").Value = -4752.2
$ws.Range("N46").Value = -8699.666999999999
$ws.Range("H93").Value = 17545922
$ws.Range("I93").Value = 23811680
$ws.Range("J93").Value = 1798.8
$ws.Range("K93").Value = 23811680
$ws.Range("L93").Value = 1798.8
$ws.Range("M93").Value = -23810432
$ws.Range("N93").Value = -4294.8
$ws.Range("H132").Value = 4198.931
$ws.Range("I132").Value = 3911.8333
$ws.Range("J132").Value = 5577
$ws.Range("K132").Value = 11735.4999
$ws.Range("L132").Value = 16731
$ws.Range("M132").Value = -9205.499899999999
$ws.Range("N132").Value = -21791

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15153808
$ws.Range("I81").Value = 20835392
$ws.Range("J81").Value = 2916.6667
$ws.Range("K81").Value = 41670784
$ws.Range("L81").Value = 5833.3334
$ws.Range("M81").Value = -41669723
$ws.Range("N81").Value = -7955.3334
$ws.Range("H84").Value = 15153808
$ws.Range("I84").Value = 20835392
$ws.Range("J84").Value = 2916.6667
$ws.Range("K84").Value = 208353920
$ws.Range("L84").Value = 29166.667
$ws.Range("M84").Value = -208348616
$ws.Range("N84").Value = -39774.667
$ws.Range("H96").Value = 5185.7144
$ws.Range("I96").Value = 4900
$ws.Range("J96").Value = 6900
$ws.Range("K96").Value = 4900
$ws.Range("L96").Value = 6900
$ws.Range("M96").Value = -3527
$ws.Range("N96").Value = -9646
$ws.Range("H122").Value = 4922
$ws.Range("I122").Value = 3960.6
$ws.Range("J122").Value = 6123.75
$ws.Range("K122").Value = 11881.8
$ws.Range("L122").Value = 18371.25
$ws.Range("M122").Value = -9431.799999999999
$ws.Range("N122").Value = -23271.25
$ws.Range("H132").Value = 25273816
$ws.Range("I132").Value = 33334424
$ws.Range("J132").Value = 1091990.8
$ws.Range("K132").Value = 100003272
$ws.Range("L132").Value = 3275972.4
$ws.Range("M132").Value = -100000742
$ws.Range("N132").Value = -3281032.4
$ws.Range("H136").Value = 2155.8525
$ws.Range("I136").Value = 1880.6875
$ws.Range("J136").Value = 3171.8462
$ws.Range("K136").Value = 5642.0625
$ws.Range("L136").Value = 9515.5386
$ws.Range("M136").Value = -3092.0625
$ws.Range("N136").Value = -14615.5386
